$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) CPPbES ("CO2 Capture Potential by Electricity Source")
#    - split "natural gas nonpeaker" into two separate technologies
#    - zero out every existing capture-rate value
#    - append new CCS / advanced generation technologies
# ---------------------------------------------------------------------------
$wsES = $wb.Worksheets.Item("CPPbES")

$wsES.Rows.Item(4).Insert()
$wsES.Range("A3").Value = "natural gas steam turbine"
$wsES.Range("B3").Value = 0
$wsES.Range("A4").Value = "natural gas combined cycle"
$wsES.Range("B4").Value = 0

$wsES.Range("B2").Value = 0
$wsES.Range("B5:B18").Value = 0

$wsES.Range("A19").Value = "hard coal w CCS"
$wsES.Range("B19").Value = 0.95
$wsES.Range("A20").Value = "natural gas combined cycle w CCS"
$wsES.Range("B20").Value = 0.95
$wsES.Range("A21").Value = "biomass w CCS"
$wsES.Range("B21").Value = 0.95
$wsES.Range("A22").Value = "lignite w CCS"
$wsES.Range("B22").Value = 0.95
$wsES.Range("A23").Value = "small modular reactor"
$wsES.Range("B23").Value = 0
$wsES.Range("A24").Value = "hydrogen combustion turbine"
$wsES.Range("B24").Value = 0
$wsES.Range("A25").Value = "hydrogen combined cycle"
$wsES.Range("B25").Value = 0

$wsES.Range("A24:A25").Font.Color = 0
$wsES.Range("A24:A25").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 2) About sheet clean-up
#    - drop the stray empty formatted cell on row 8
#    - clear now-redundant formatting left over from an old highlighted block
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B8").Clear()
$wsAbout.Range("B4").ClearFormats()
$wsAbout.Range("A10").ClearFormats()
$wsAbout.Range("D17:E21").ClearFormats()

# ---------------------------------------------------------------------------
# 3) New sheet: CPPbHS ("CO2 Capture Potential by Hydrogen Source")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHS = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsHS.Name = "CPPbHS"

$wsHS.Tab.ColorIndex = $wb.Worksheets.Item("CPPbI").Tab.ColorIndex

$wsHS.Columns.Item(1).ColumnWidth = 45.42578125
$wsHS.Columns.Item(2).ColumnWidth = 24.7109375
$wsHS.Columns.Item(3).ColumnWidth = 25.85546875

$wsHS.Range("A1").Value = "Unit: dimentionless (fraction of CO2 capturable)"
$wsHS.Range("A1").Font.Italic = $true
$wsHS.Range("B1").Value = "capture rate"
$wsHS.Range("A2").Value = "natural gas reforming with CCS"
$wsHS.Range("B2").Value = 0.85

$wsAbout.Activate()
$wsAbout.Range("A1").Select()

Write-Host "edit complete"
